$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12. This shifts the existing rows 12..132
# down to 13..133, which matches the target workbook's data (each old
# row's content reappears one row further down, and a fresh row 133 is
# created holding what used to be row 132).
$ws.Rows.Item(12).Insert()

# Populate the newly-inserted row 12 with this week's record: the same
# static reference values used throughout the sheet, plus the new date
# and the same price/volume tuple as the most recent prior entry for
# "Red Globe" / "Primera" at "Región de O'Higgins".
$ws.Cells.Item(12, 1).Value = 7
$ws.Cells.Item(12, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(12, 3).Value = "Ñuble"
$ws.Cells.Item(12, 4).Value = 44649
$ws.Cells.Item(12, 5).Value = 16
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100109
$ws.Cells.Item(12, 8).Value = "Uva"
$ws.Cells.Item(12, 9).Value = 100109001
$ws.Cells.Item(12, 10).Value = "Uva"
$ws.Cells.Item(12, 11).Value = "Red Globe"
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 120
$ws.Cells.Item(12, 14).Value = 9000
$ws.Cells.Item(12, 15).Value = 10000
$ws.Cells.Item(12, 16).Value = 9500
$ws.Cells.Item(12, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(12, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(12, 19).Value = 528
$ws.Cells.Item(12, 20).Value = 18
